# agregar contador de documentos a incidencias
#
# This script rewrites the data rows of the "AIO" and "Otros" sheets so
# that each row gets an extra "document counter" style duplicate entry
# (rows 5-6 on AIO shift to 7-8, rows 4 on Otros shifts to 4-6) and every
# data cell (including what used to be numeric N-Serie / Codigo Inventario
# columns) is stored as plain text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AIO": grows from A1:J6 to A1:J8
# ---------------------------------------------------------------------
$wsAio = $wb.Worksheets.Item("AIO")

$aioData = @(
  @('Concepcion', 'Coronel', 'OFICINA', '01374859', 'BALDOMERO LILLO', 'AIO', 'Lenovo', 'M700z', '0', '8013913'),
  @('Concepcion', 'Coronel', 'OFICINA', '02993344', 'BALDOMERO LILLO', 'AIO', 'Lenovo', '71z', '0', '8013913'),
  @('Concepcion', 'Concepcion', 'CLASICO', '01773356', 'ABKELAY KIMUN', 'AIO', 'Lenovo', 'E73z', '0', '8013913'),
  @('V530', '2', 'OFICINA', '01233421', 'BALDOMERO LILLO', 'AIO', 'Lenovo', 'Coronel', '0', '8013913'),
  @('V510z', '1', $null, '02988776', 'VILLA GENESIS', 'AIO', 'Lenovo', 'Los Angeles', '0', '17006716'),
  @('V530', '2', 'OFICINA', '01233421', 'BALDOMERO LILLO', 'AIO', 'Lenovo', 'Coronel', '0', '8013913'),
  @('V510z', '1', $null, '02988776', 'VILLA GENESIS', 'AIO', 'Lenovo', 'Los Angeles', '0', '17006716')
)

$wsAio.Range("A2:J8").NumberFormat = "@"

for ($i = 0; $i -lt $aioData.Length; $i++) {
  $row = $i + 2
  for ($j = 0; $j -lt 10; $j++) {
    $col = $j + 1
    $val = $aioData[$i][$j]
    if ($null -eq $val) {
      $wsAio.Cells.Item($row, $col).ClearContents()
    } else {
      $wsAio.Cells.Item($row, $col).Value = $val
    }
  }
}

# ---------------------------------------------------------------------
# Sheet "Otros": grows from A1:J4 to A1:J6
# ---------------------------------------------------------------------
$wsOtros = $wb.Worksheets.Item("Otros")

$otrosData = @(
  @('Concepcion', 'Coronel', 'OFICINA', '001', 'BALDOMERO LILLO', 'Mouse', 'ACER', 'S3 SERIES', '0001', '0001'),
  @('Concepcion', 'Concepcion', 'CLASICO', '123', 'ABKELAY KIMUN', 'Teclado', 'ACER', 'S3 SERIES', '123', '123'),
  @('Concepcion', 'Coronel', 'OFICINA', '00137', 'BALDOMERO LILLO', 'Teclado', 'Lenovo', 'V510z', '137', '137'),
  @('Biobío', 'Los Angeles', $null, '2A', 'VILLA GENESIS', 'CPU', 'HP', '340', '2', '2'),
  @('Concepcion', 'Santa Juana', $null, '356734', 'Nueva-Unidad', 'Mouse', 'ACER', 'S3 SERIES', '2346', '0839433')
)

$wsOtros.Range("A2:J6").NumberFormat = "@"

for ($i = 0; $i -lt $otrosData.Length; $i++) {
  $row = $i + 2
  for ($j = 0; $j -lt 10; $j++) {
    $col = $j + 1
    $val = $otrosData[$i][$j]
    if ($null -eq $val) {
      $wsOtros.Cells.Item($row, $col).ClearContents()
    } else {
      $wsOtros.Cells.Item($row, $col).Value = $val
    }
  }
}
